$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to text formatting before writing so that
# Excel does not auto-coerce numeric-looking strings (e.g. "0.9998") into
# actual numbers, which would lose the original text representation.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.123.23'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.824.98'
$ws.Range("E3").Value = '  -0.53%  '
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '241.57'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").Value = '0.6159'
$ws.Range("E6").Value = '  -1.90%  '
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  -0.20%  '
$ws.Range("D8").Value = '0.07334'
$ws.Range("E8").Value = '  -2.43%  '
$ws.Range("D9").Value = '0.2894'
$ws.Range("E9").Value = '  -1.07%  '
$ws.Range("D10").Value = '22.92'
$ws.Range("E10").Value = '  -1.20%  '
$ws.Range("D11").Value = '0.07661'
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = '1.814.57'
$ws.Range("E12").Value = '  -1.10%  '
$ws.Range("D13").Value = '4.952'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '0.6612'
$ws.Range("E14").Value = '  -1.06%  '
$ws.Range("D15").Value = '81.71'
$ws.Range("E15").Value = '  -1.14%  '
$ws.Range("D16").Value = '0.000008945'
$ws.Range("E16").Value = '  -4.54%  '
$ws.Range("D17").Value = '5.855'
$ws.Range("E17").Value = '  -2.36%  '
$ws.Range("D18").Value = '29.080.76'
$ws.Range("E18").Value = '  -0.11%  '
$ws.Range("D19").Value = '2.049.54'
$ws.Range("E19").Value = '  -1.39%  '
$ws.Range("D20").Value = '237.66'
$ws.Range("E20").Value = '  +6.28%  '
$ws.Range("D21").Value = '12.47'
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D22").Value = '0.9995'
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("D23").Value = '7.134'
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").Value = '158.35'
$ws.Range("D26").Value = '0.1409'
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("D27").Value = '8.431'
$ws.Range("E27").Value = '  -0.85%  '
$ws.Range("D28").Value = '17.63'
$ws.Range("E28").Value = '  -1.63%  '
$ws.Range("E29").Value = '  -1.20%  '
$ws.Range("D30").Value = '0.05585'
$ws.Range("E30").Value = '  -1.82%  '
$ws.Range("D31").Value = '4.094'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("D32").Value = '4.100'
$ws.Range("E32").Value = '  -1.35%  '
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").Value = '1.826'
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("D35").Value = '0.7338'
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("D36").Value = '1.131'
$ws.Range("E36").Value = '  -0.76%  '
$ws.Range("D37").Value = '2.622'
$ws.Range("E37").Value = '  -1.85%  '
$ws.Range("D38").Value = '2.827'
$ws.Range("E38").Value = '  +2.31%  '
$ws.Range("D39").Value = '1.204.00'
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("D40").Value = '0.01757'
$ws.Range("E40").Value = '  -1.45%  '
$ws.Range("D41").Value = '6.379'
$ws.Range("E41").Value = '  -2.19%  '
$ws.Range("D42").Value = '0.8941'
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").Value = '0.9994'
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").Value = '100.68'
$ws.Range("E44").Value = '  -1.24%  '
$ws.Range("D45").Value = '1.957.85'
$ws.Range("E45").Value = '  -0.97%  '
$ws.Range("D46").Value = '64.60'
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("E47").Value = '  -2.71%  '
$ws.Range("D48").Value = '0.5074'
$ws.Range("E48").Value = '  -0.32%  '
$ws.Range("D49").Value = '9.069'
$ws.Range("E49").Value = '  +0.67%  '
$ws.Range("D50").Value = '0.3993'
$ws.Range("E50").Value = '  -2.13%  '
$ws.Range("D51").Value = '0.05796'
$ws.Range("E51").Value = '  -0.51%  '

# Restore the cells to their original (unstyled) state now that the text
# values are committed, so no stray number-format styling is left behind.
$dataRange.Style = "Normal"

